# Workbook-level changes: rename sheets, tweak the saved window view.
$wb = $excel.ActiveWorkbook

$wsCodes = $wb.Worksheets.Item(1)   # "Error codes"      -> "Message codes"
$wsTypes = $wb.Worksheets.Item(2)   # "Error types"       -> "Codes descriptions"
$wsTests = $wb.Worksheets.Item(3)   # "Error tests"       (unchanged)

$wsCodes.Name = "Message codes"
$wsTypes.Name = "Codes descriptions"

# Saved window position/size (best effort - mirrors the workbook window move/resize).
$win = $wb.Windows.Item(1)
$win.Top = 165
$win.Height = 6180

# --- Sheet "Message codes" (sheet1) -----------------------------------------
# Fix a long-standing column swap on row 51 (Error / Error code were reversed).
$oldA51 = $wsCodes.Range("A51").Value2
$oldB51 = $wsCodes.Range("B51").Value2
$wsCodes.Range("A51").Value2 = $oldB51
$wsCodes.Range("B51").Value2 = $oldA51

# New row 41: junit tests for the refresh-status "unknown cause" ack error.
$wsCodes.Range("A41").Value2 = "Ack is KO, and no errors were found in it => unknown cause"
$wsCodes.Range("B41").Value2 = "ERR806"
$wsCodes.Range("C41").Value2 = "yes"

# New row 74: junit tests for the refresh-status success code.
$wsCodes.Range("A74").Value2 = "Refresh status successfully completed"
$wsCodes.Range("B74").Value2 = "OK500"
$wsCodes.Range("C74").Value2 = "yes"

# Resize the table to cover the newly added rows (table grows a bit further
# than the populated data, matching the author's manual resize).
$tblCodes = $wsCodes.ListObjects.Item(1)
$tblCodes.Resize($wsCodes.Range("A1:C78"))

# Move the selection to reflect where the author ended up editing.
$wsCodes.Range("A42").Select()

# --- Sheet "Codes descriptions" (sheet2) ------------------------------------
$wsTypes.Range("A1").Value2 = "Code group"
$wsTypes.Range("B1").Value2 = "Description"

$wsTypes.Range("A6").Select()

# "Message codes" stays the tab that is active/selected when the file is reopened.
$wsCodes.Activate()
$wsCodes.Range("A42").Select()
